$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns C through H (data_path_2, data_path_3, index_0..index_3)
$ws.Columns("C:H").Delete()

# Update the remaining header values
$ws.Range("A1").Value = "input_useremail"
$ws.Range("B1").Value = "input_useremail_1"

# Clear the data row values so A2/B2 become blank
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()

# Set the new column widths
$ws.Columns.Item(1).ColumnWidth = 16.1
$ws.Columns.Item(2).ColumnWidth = 18.15
